{"js": "// 1. \"Once the signed agreement ...\" paragraph:\n//    - \"Google drive\" -> \"One drive\"\n//    - \"Dr. Liang Lin\" -> \"Prof. Liang Lin\"\n//    - The Word \"last edit position\" bookmark (_GoBack) ends up right before\n//      \"drive\" in \"One drive\" (after the word \"One\" and the following space).\nconst body = context.document.body;\n\n// --- Replace \"Google drive\" with \"One drive\" ---\nconst googleDrive = body.search(\"Google drive\", { matchCase: true, matchWholeWord: false });\ngoogleDrive.load(\"text\");\nawait context.sync();\nif (googleDrive.items.length > 0) {\n  googleDrive.items[0].insertText(\"One drive\", Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// --- Replace \"Dr. Liang Lin\" with \"Prof. Liang Lin\" (inside same paragraph) ---\nconst drLiangLin = body.search(\"Dr. Liang Lin\", { matchCase: true });\ndrLiangLin.load(\"text\");\nawait context.sync();\nif (drLiangLin.items.length > 0) {\n  drLiangLin.items[0].insertText(\"Prof. Liang Lin\", Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// --- Move the _GoBack bookmark to right before \"drive to download\" ---\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\nconst driveToDownload = body.search(\"drive to download\", { matchCase: true });\nawait context.sync();\nif (driveToDownload.items.length > 0) {\n  const insertionPoint = driveToDownload.items[0].getRange(\"Start\");\n  insertionPoint.insertBookmark(\"_GoBack\");\n  await context.sync();\n}\n\n// 2. Bibliography entry 1 (Yuan Xie et al.): join the two runs that used to be\n//    split by the (now relocated) _GoBack bookmark, and tighten \"MM '20).\" to\n//    \"MM'20).\"\nconst mmSpace = body.search(\"MM '20).\", { matchCase: true });\nmmSpace.load(\"text\");\nawait context.sync();\nif (mmSpace.items.length > 0) {\n  mmSpace.items[0].insertText(\"MM'20).\", Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// 3. Bibliography entry 2 (Tianshui Chen et al.): normalize the all-caps\n//    journal title to title case.\nconst tpami = body.search(\n  \"In IEEE TRANSACTIONS ON PATTERN ANALYSIS AND MACHINE INTELLIGENCE (TPAMI\",\n  { matchCase: true }\n);\ntpami.load(\"text\");\nawait context.sync();\nif (tpami.items.length > 0) {\n  tpami.items[0].insertText(\n    \"In IEEE Transactions on Pattern Analysis and Machine Intelligence (TPAMI\",\n    Word.InsertLocation.replace\n  );\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\nfunction Replace-FirstMatch($doc, $searchText, $replaceText) {\n    $find = $doc.Content.Find\n    $find.ClearFormatting()\n    $find.Text = $searchText\n    $found = $find.Execute()\n    if ($found) {\n        # Assigning .Text directly (instead of using Find.Replacement) keeps\n        # the run's original formatting and avoids Word's Find/Replace\n        # \"smart quotes\" auto-correction of straight apostrophes.\n        $find.Parent.Text = $replaceText\n    }\n    return $found\n}\n\n# --- 1a. \"Google drive\" -> \"One drive\" ---\nReplace-FirstMatch $d \"Google drive\" \"One drive\" | Out-Null\n\n# --- 1b. \"Dr. Liang Lin\" -> \"Prof. Liang Lin\" ---\nReplace-FirstMatch $d \"Dr. Liang Lin\" \"Prof. Liang Lin\" | Out-Null\n\n# --- 1c. Move Word's \"last edit\" (_GoBack) bookmark from the bibliography\n#         paragraph (where it originally sat) to right before\n#         \"drive to download\" in the paragraph we just edited. ---\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks.Item(\"_GoBack\").Delete()\n}\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Text = \"drive to download\"\n$find.Execute() | Out-Null\n$insertionPoint = $d.Range($find.Parent.Start, $find.Parent.Start)\n$d.Bookmarks.Add(\"_GoBack\", $insertionPoint) | Out-Null\n\n# --- 2. \"MM '20).\" -> \"MM'20).\" ---\nReplace-FirstMatch $d \"MM '20).\" \"MM'20).\" | Out-Null\n\n# --- 3. Normalize the all-caps TPAMI journal title to title case. ---\nReplace-FirstMatch $d \"In IEEE TRANSACTIONS ON PATTERN ANALYSIS AND MACHINE INTELLIGENCE (TPAMI\" \"In IEEE Transactions on Pattern Analysis and Machine Intelligence (TPAMI\" | Out-Null\n"}
